# Agregar contenido clase 9 y control 3
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clase 9 (row 11): move the "material" text from the "material_futuro" column (I)
# into the real "material" column (G), now that the slides are available.
$ws.Range("G11").Value2 = $ws.Range("I11").Value2
$ws.Range("I11").ClearContents()

# Control 3 (row 12): move the "entrega" text from the "entrega_futuro" column (H)
# into the real "entrega" column (F), now that the assignment is due.
$ws.Range("F12").Value2 = $ws.Range("H12").Value2
$ws.Range("H12").ClearContents()

# Update the saved selection to match the author's final cursor position.
$ws.Range("F12").Select()
